$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @(0.21,0.22,0.5600000000000001),
  @(0.15,0.68,0.17),
  @(0.65,0.18,0.17),
  @(0.71,0.14,0.14),
  @(0.23,0.53,0.24),
  @(0.16,0.72,0.13),
  @(0.19,0.23,0.58),
  @(0.18,0.63,0.19),
  @(0.17,0.65,0.18),
  @(0.3,0.25,0.45),
  @(0.17,0.5600000000000001,0.27),
  @(0.18,0.62,0.19),
  @(0.12,0.75,0.12),
  @(0.73,0.13,0.14),
  @(0.17,0.23,0.6),
  @(0.19,0.18,0.63),
  @(0.79,0.1,0.1),
  @(0.68,0.16,0.16),
  @(0.23,0.57,0.21),
  @(0.75,0.11,0.14),
  @(0.16,0.19,0.66),
  @(0.2,0.54,0.25),
  @(0.15,0.71,0.14),
  @(0.79,0.1,0.11),
  @(0.75,0.13,0.12),
  @(0.74,0.13,0.13),
  @(0.18,0.18,0.64),
  @(0.12,0.11,0.77),
  @(0.2,0.54,0.26),
  @(0.13,0.75,0.13),
  @(0.13,0.7,0.17),
  @(0.64,0.17,0.19),
  @(0.13,0.72,0.15),
  @(0.79,0.1,0.11),
  @(0.21,0.21,0.58),
  @(0.64,0.18,0.18),
  @(0.18,0.62,0.19),
  @(0.2,0.5,0.3),
  @(0.7,0.16,0.14),
  @(0.17,0.67,0.16),
  @(0.63,0.18,0.18),
  @(0.22,0.59,0.19),
  @(0.8,0.1,0.1),
  @(0.17,0.14,0.68),
  @(0.68,0.16,0.16),
  @(0.2,0.6,0.2),
  @(0.17,0.65,0.18),
  @(0.18,0.6,0.22),
  @(0.23,0.55,0.22),
  @(0.2,0.21,0.59),
  @(0.22,0.17,0.6),
  @(0.13,0.13,0.74),
  @(0.53,0.25,0.22),
  @(0.12,0.13,0.76),
  @(0.2,0.2,0.6),
  @(0.78,0.11,0.11),
  @(0.14,0.16,0.7),
  @(0.1,0.11,0.79),
  @(0.36,0.41,0.24),
  @(0.84,0.07000000000000001,0.09),
  @(0.38,0.37,0.25),
  @(0.6899999999999999,0.16,0.16),
  @(0.82,0.09,0.09),
  @(0.2,0.51,0.29),
  @(0.12,0.13,0.75),
  @(0.22,0.55,0.23),
  @(0.57,0.25,0.18),
  @(0.17,0.64,0.19),
  @(0.13,0.12,0.75),
  @(0.25,0.25,0.51),
  @(0.78,0.11,0.11),
  @(0.17,0.58,0.25),
  @(0.15,0.68,0.16),
  @(0.14,0.16,0.7),
  @(0.83,0.09,0.09),
  @(0.2,0.21,0.59),
  @(0.17,0.16,0.68),
  @(0.74,0.13,0.13),
  @(0.09,0.11,0.8),
  @(0.19,0.5,0.31),
  @(0.23,0.53,0.24),
  @(0.13,0.16,0.71),
  @(0.18,0.63,0.19),
  @(0.76,0.12,0.13),
  @(0.13,0.17,0.7),
  @(0.39,0.27,0.34),
  @(0.67,0.16,0.18),
  @(0.53,0.23,0.23),
  @(0.19,0.22,0.58),
  @(0.77,0.12,0.11)
)

for ($i = 0; $i -lt $data.Count; $i++) {
  $row = $data[$i]
  $ws.Cells.Item($i + 2, 2).Value = $row[0]
  $ws.Cells.Item($i + 2, 3).Value = $row[1]
  $ws.Cells.Item($i + 2, 4).Value = $row[2]
}
